$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-25 Monday" "2024-11-26 Tuesday"

Replace-Text "405÷5=" "711÷6="
Replace-Text "730÷4=" "174÷7="
Replace-Text "863÷3=" "755÷3="
Replace-Text "632÷4=" "647÷6="
Replace-Text "633÷2=" "333÷8="

Replace-Text "239÷4=" "599÷5="
Replace-Text "258÷2=" "847÷3="
Replace-Text "627÷7=" "978÷2="
Replace-Text "841÷4=" "840÷5="
Replace-Text "927÷3=" "563÷8="

Replace-Text "165÷9=" "977÷3="
Replace-Text "146÷6=" "766÷8="
Replace-Text "488÷5=" "351÷5="
Replace-Text "997÷5=" "292÷3="
Replace-Text "794÷2=" "754÷6="

Replace-Text "744÷6=" "758÷4="
Replace-Text "123÷6=" "976÷9="
Replace-Text "997÷8=" "629÷7="
Replace-Text "912÷6=" "410÷9="
Replace-Text "146÷9=" "530÷7="

Replace-Text "651÷8=" "496÷4="
Replace-Text "883÷7=" "442÷5="
Replace-Text "957÷4=" "779÷8="
Replace-Text "618÷9=" "840÷9="
Replace-Text "582÷2=" "489÷8="
